$wb = $excel.ActiveWorkbook

# --- 1. Remove stray empty cells from the "ODI Batting" sheet (B2, B8, B14) ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B8").ClearContents()
$batting.Range("B14").ClearContents()

# --- 2. Add the new "ODI Batting Extra" sheet after "ODI Bowling" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# --- 3. Header row, styled like the other sheets header rows ---
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"
$batting.Range("A1:F1").Copy()
$newSheet.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4. Data rows. Columns A, C, D, E hold numeric-looking text, so format as Text first
#        to keep them as strings rather than auto-converted numbers/percentages. ---
$newSheet.Range("A2:A20").NumberFormat = "@"
$newSheet.Range("C2:E20").NumberFormat = "@"

$newSheet.Range("A2").Value = "4563"
$newSheet.Range("B2").Value = 6
$newSheet.Range("F2").Value = "NO"

$newSheet.Range("A3").Value = "4566"
$newSheet.Range("B3").Value = 6
$newSheet.Range("C3").Value = "0"
$newSheet.Range("D3").Value = "0"
$newSheet.Range("E3").Value = "0.38%"
$newSheet.Range("F3").Value = "NO"

$newSheet.Range("A4").Value = "4568"
$newSheet.Range("F4").Value = "NO"

$newSheet.Range("A5").Value = "4605"
$newSheet.Range("B5").Value = 7
$newSheet.Range("C5").Value = "10"
$newSheet.Range("D5").Value = "7"
$newSheet.Range("E5").Value = "41.64%"
$newSheet.Range("F5").Value = "YES"

$newSheet.Range("A6").Value = "4608"
$newSheet.Range("B6").Value = 7
$newSheet.Range("C6").Value = "3"
$newSheet.Range("D6").Value = "3"
$newSheet.Range("E6").Value = "19.18%"
$newSheet.Range("F6").Value = "YES"

$newSheet.Range("A7").Value = "4614"
$newSheet.Range("B7").Value = 7
$newSheet.Range("C7").Value = "1"
$newSheet.Range("D7").Value = "1"
$newSheet.Range("E7").Value = "5.83%"
$newSheet.Range("F7").Value = "NO"

$newSheet.Range("A8").Value = "4625"
$newSheet.Range("B8").Value = 7
$newSheet.Range("F8").Value = "NO"

$newSheet.Range("A9").Value = "4636"
$newSheet.Range("B9").Value = 7
$newSheet.Range("C9").Value = "5"
$newSheet.Range("D9").Value = "0"
$newSheet.Range("E9").Value = "16.32%"
$newSheet.Range("F9").Value = "NO"

$newSheet.Range("A10").Value = "4639"
$newSheet.Range("B10").Value = 6
$newSheet.Range("C10").Value = "0"
$newSheet.Range("D10").Value = "0"
$newSheet.Range("E10").Value = "2.83%"
$newSheet.Range("F10").Value = "NO"

$newSheet.Range("A11").Value = "4642"
$newSheet.Range("F11").Value = "NO"

$newSheet.Range("A12").Value = "4647"
$newSheet.Range("F12").Value = "NO"

$newSheet.Range("A13").Value = "4648"
$newSheet.Range("B13").Value = 6
$newSheet.Range("C13").Value = "1"
$newSheet.Range("D13").Value = "0"
$newSheet.Range("E13").Value = "14.63%"
$newSheet.Range("F13").Value = "NO"

$newSheet.Range("A14").Value = "4673"
$newSheet.Range("F14").Value = "NO"

$newSheet.Range("A15").Value = "4686"
$newSheet.Range("F15").Value = "NO"

$newSheet.Range("A16").Value = "4688"
$newSheet.Range("B16").Value = 7
$newSheet.Range("C16").Value = "0"
$newSheet.Range("D16").Value = "0"
$newSheet.Range("E16").Value = "3.07%"
$newSheet.Range("F16").Value = "NO"

$newSheet.Range("A17").Value = "4690"
$newSheet.Range("F17").Value = "NO"

$newSheet.Range("A18").Value = "4692"
$newSheet.Range("F18").Value = "NO"

$newSheet.Range("A19").Value = "4695"
$newSheet.Range("B19").Value = 7
$newSheet.Range("C19").Value = "4"
$newSheet.Range("D19").Value = "0"
$newSheet.Range("E19").Value = "20.37%"
$newSheet.Range("F19").Value = "NO"

$newSheet.Range("A20").Value = "4697"
$newSheet.Range("B20").Value = 7
$newSheet.Range("C20").Value = "3"
$newSheet.Range("D20").Value = "1"
$newSheet.Range("E20").Value = "8.81%"
$newSheet.Range("F20").Value = "NO"

# --- 5. Restore the original active sheet selection ---
$wb.Worksheets.Item("Player Info").Activate()